$p = $ppt.ActivePresentation

# Slide 4 ("ARBORESCENCES"): move the "ZoneTexte 9" body textbox up slightly
# (title/body repositioning - changer la position du titre).
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(4)
$shape4.Top = 204.14394380787402

# Slide 5 ("EXPLICATIONS DES PARTAGES"): reposition the "ZoneTexte 9" ("TYPE DE DROIT") textbox.
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(3)
$shape5.Left = 20.611338682677168
$shape5.Top = 167.94141732283464
